# Auto-generated edit script
# Applies the cell-value changes described by the commit diff across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets of the Famfrit profits workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 839.5
$ws.Range("J53").Value = 2156.8
$ws.Range("L53").Value = 2156.8
$ws.Range("N53").Value = -3430.8
$ws.Range("H64").Value = 10699.4
$ws.Range("I64").Value = 10374.25
$ws.Range("K64").Value = 10374.25
$ws.Range("M64").Value = -10126.25
$ws.Range("H67").Value = 10699.4
$ws.Range("I67").Value = 10374.25
$ws.Range("K67").Value = 10374.25
$ws.Range("M67").Value = -9516.25
$ws.Range("H70").Value = 2322.52
$ws.Range("I70").Value = 2010.6364
$ws.Range("J70").Value = 2567.5715
$ws.Range("K70").Value = 6031.9092
$ws.Range("L70").Value = 7702.7145
$ws.Range("M70").Value = -5761.9092
$ws.Range("N70").Value = -8242.7145
$ws.Range("H73").Value = 2322.52
$ws.Range("I73").Value = 2010.6364
$ws.Range("J73").Value = 2567.5715
$ws.Range("K73").Value = 6031.9092
$ws.Range("L73").Value = 7702.7145
$ws.Range("M73").Value = -5095.9092
$ws.Range("N73").Value = -9574.7145
$ws.Range("H118").Value = 1598.9231
$ws.Range("I118").Value = 1617.6364
$ws.Range("K118").Value = 4852.9092
$ws.Range("M118").Value = -3195.9092
$ws.Range("H137").Value = 2440.625
$ws.Range("I137").Value = 2451.3684
$ws.Range("J137").Value = 2399.8
$ws.Range("K137").Value = 7354.1052
$ws.Range("L137").Value = 7199.400000000001
$ws.Range("M137").Value = -4804.1052
$ws.Range("N137").Value = -12299.4
$ws.Range("H138").Value = 6671253.5
$ws.Range("I138").Value = 1734.7646
$ws.Range("K138").Value = 5204.293799999999
$ws.Range("M138").Value = -64.29379999999946

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14936743
$ws.Range("I32").Value = 23814706
$ws.Range("J32").Value = 21767.72
$ws.Range("K32").Value = 23814706
$ws.Range("L32").Value = 21767.72
$ws.Range("M32").Value = -23814419
$ws.Range("N32").Value = -22341.72
$ws.Range("H74").Value = 40047044
$ws.Range("I74").Value = 62571372
$ws.Range("J74").Value = 3788.6667
$ws.Range("K74").Value = 62571372
$ws.Range("L74").Value = 3788.6667
$ws.Range("M74").Value = -62570498
$ws.Range("N74").Value = -5536.6667
$ws.Range("H77").Value = 40047044
$ws.Range("I77").Value = 62571372
$ws.Range("J77").Value = 3788.6667
$ws.Range("K77").Value = 312856860
$ws.Range("L77").Value = 18943.3335
$ws.Range("M77").Value = -312852492
$ws.Range("N77").Value = -27679.3335
$ws.Range("H132").Value = 24396742
$ws.Range("I132").Value = 7490
$ws.Range("K132").Value = 22470
$ws.Range("M132").Value = -19940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2103.3125
$ws.Range("I94").Value = 1812.28
$ws.Range("K94").Value = 1812.28
$ws.Range("M94").Value = -1361.28

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17862960
$ws.Range("I31").Value = 5120.381
$ws.Range("K31").Value = 5120.381
$ws.Range("M31").Value = -4825.381
$ws.Range("H34").Value = 17862960
$ws.Range("I34").Value = 5120.381
$ws.Range("K34").Value = 5120.381
$ws.Range("M34").Value = -4918.381
$ws.Range("H62").Value = 3664.889
$ws.Range("I62").Value = 2499.8333
$ws.Range("J62").Value = 5995
$ws.Range("K62").Value = 2499.8333
$ws.Range("L62").Value = 5995
$ws.Range("M62").Value = -1875.8333
$ws.Range("N62").Value = -7243
$ws.Range("H65").Value = 3664.889
$ws.Range("I65").Value = 2499.8333
$ws.Range("J65").Value = 5995
$ws.Range("K65").Value = 12499.1665
$ws.Range("L65").Value = 29975
$ws.Range("M65").Value = -9379.166499999999
$ws.Range("N65").Value = -36215
$ws.Range("H122").Value = 1644.2941
$ws.Range("I122").Value = 1637.8334
$ws.Range("J122").Value = 1659.8
$ws.Range("K122").Value = 4913.5002
$ws.Range("L122").Value = 4979.4
$ws.Range("M122").Value = -2463.5002
$ws.Range("N122").Value = -9879.4
$ws.Range("H132").Value = 55144.848
$ws.Range("I132").Value = 73069.34
$ws.Range("K132").Value = 219208.02
$ws.Range("M132").Value = -216678.02
$ws.Range("H134").Value = 1818
$ws.Range("I134").Value = 1619.6364
$ws.Range("K134").Value = 4858.9092
$ws.Range("M134").Value = -2323.9092
$ws.Range("H141").Value = 287243.12
$ws.Range("J141").Value = 308148.5
$ws.Range("L141").Value = 308148.5
$ws.Range("N141").Value = -318508.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 162497.38
$ws.Range("J37").Value = 162497.38
$ws.Range("L37").Value = 487492.14
$ws.Range("N37").Value = -487716.14
$ws.Range("H122").Value = 1539.2222
$ws.Range("I122").Value = 386.75
$ws.Range("K122").Value = 3480.75
$ws.Range("M122").Value = -1030.75
$ws.Range("H131").Value = 35752.79
$ws.Range("I131").Value = 105020
$ws.Range("J131").Value = 5636.609
$ws.Range("K131").Value = 315060
$ws.Range("L131").Value = 16909.827
$ws.Range("M131").Value = -310020
$ws.Range("N131").Value = -26989.827
$ws.Range("H137").Value = 2607.6553
$ws.Range("J137").Value = 3035.5
$ws.Range("L137").Value = 9106.5
$ws.Range("N137").Value = -19306.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13599.25
$ws.Range("I43").Value = 2700
$ws.Range("J43").Value = 24498.5
$ws.Range("K43").Value = 2700
$ws.Range("L43").Value = 24498.5
$ws.Range("M43").Value = -2549
$ws.Range("H80").Value = 3315.7856
$ws.Range("I80").Value = 3010.875
$ws.Range("J80").Value = 3722.3333
$ws.Range("K80").Value = 3010.875
$ws.Range("L80").Value = 3722.3333
$ws.Range("M80").Value = -2012.875
$ws.Range("N80").Value = -5718.3333
$ws.Range("H83").Value = 3315.7856
$ws.Range("I83").Value = 3010.875
$ws.Range("J83").Value = 3722.3333
$ws.Range("K83").Value = 15054.375
$ws.Range("L83").Value = 18611.6665
$ws.Range("M83").Value = -10062.375
$ws.Range("N83").Value = -28595.6665
$ws.Range("H132").Value = 2618.6667
$ws.Range("I132").Value = 2541.2632
$ws.Range("K132").Value = 7623.7896
$ws.Range("M132").Value = -5093.7896
$ws.Range("N43").Value = -24800.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5682.625
$ws.Range("I40").Value = 4960.8
$ws.Range("J40").Value = 6885.6665
$ws.Range("K40").Value = 4960.8
$ws.Range("L40").Value = 6885.6665
$ws.Range("M40").Value = -4824.8
$ws.Range("N40").Value = -7157.6665
$ws.Range("H46").Value = 1575.9246
$ws.Range("I46").Value = 950.2105
$ws.Range("J46").Value = 3161.0667
$ws.Range("K46").Value = 950.2105
$ws.Range("L46").Value = 3161.0667
$ws.Range("M46").Value = -762.2105
$ws.Range("N46").Value = -3537.0667
$ws.Range("H55").Value = 655.1429000000001
$ws.Range("I55").Value = 347.66666
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 347.66666
$ws.Range("L55").Value = 2500
$ws.Range("M55").Value = -174.66666
$ws.Range("N55").Value = -2846
$ws.Range("H64").Value = 17856.666
$ws.Range("J64").Value = 17856.666
$ws.Range("L64").Value = 17856.666
$ws.Range("N64").Value = -18306.666
$ws.Range("H67").Value = 17856.666
$ws.Range("J67").Value = 17856.666
$ws.Range("L67").Value = 17856.666
$ws.Range("N67").Value = -19416.666
$ws.Range("H122").Value = 4119.6772
$ws.Range("I122").Value = 3645.682
$ws.Range("J122").Value = 5278.3335
$ws.Range("K122").Value = 10937.046
$ws.Range("L122").Value = 15835.0005
$ws.Range("M122").Value = -8487.045999999998
$ws.Range("N122").Value = -20735.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7999.125
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 7999.125
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 7999.125
$ws.Range("N62").Value = -9247.125
$ws.Range("H63").Value = 15500
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H65").Value = 7999.125
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 7999.125
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 39995.625
$ws.Range("N65").Value = -46235.625
$ws.Range("H66").Value = 15500
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H132").Value = 5938.108
$ws.Range("I132").Value = 5908.8887
$ws.Range("J132").Value = 6990
$ws.Range("K132").Value = 17726.6661
$ws.Range("L132").Value = 20970
$ws.Range("M132").Value = -15196.6661
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

Write-Host "Applied 216 value updates, 1 new cells, 2 cleared cells."
